# Update gh-pages output data (想去人数 / F column counts, and one date value)
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 42360
$ws1.Range("F4").Value  = 10231
$ws1.Range("F5").Value  = 237
$ws1.Range("F6").Value  = 1101
$ws1.Range("F14").Value = 819
$ws1.Range("F15").Value = 354
$ws1.Range("F16").Value = 1691
$ws1.Range("F18").Value = 852
$ws1.Range("F21").Value = 738
$ws1.Range("F22").Value = 832
$ws1.Range("F23").Value = 41
$ws1.Range("F26").Value = 583
$ws1.Range("F35").Value = 193
$ws1.Range("F36").Value = 511
$ws1.Range("F37").Value = 1513
$ws1.Range("F39").Value = 1345
$ws1.Range("F46").Value = 26

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 168
$ws2.Range("F7").Value = 68

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 501

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 10231
$ws4.Range("F5").Value  = 1101
$ws4.Range("F7").Value  = 501
$ws4.Range("F10").Value = 168
$ws4.Range("F14").Value = 68
$ws4.Range("F15").Value = 819
$ws4.Range("F16").Value = 354
$ws4.Range("F17").Value = 1691
$ws4.Range("F19").Value = 852
$ws4.Range("F22").Value = 738
$ws4.Range("F23").Value = 832
$ws4.Range("F24").Value = 41
$ws4.Range("F27").Value = 583
$ws4.Range("F38").Value = 193
$ws4.Range("F42").Value = 1345
